# Trade #122 closed at 2026-02-16 21:45:17 - leadlag DOWN +0.000%
#
# This script:
#  1. Updates the rolled-up OVERALL / leadlag STRATEGY stats on the
#     "Summary" sheet (two new closed trades change trade count / win
#     rate / total P&L%).
#  2. Updates the "leadlag" sheet: trade #100 (row 76) and #101 (row 77)
#     which were still OPEN are now CLOSED with exit price / P&L filled
#     in, and a brand-new trade #122 (row 97) is appended as OPEN.
#  3. Mirrors the two newly-closed trades (#100 / #101) onto the
#     "All Trades" sheet as new rows 101 / 102.
#  4. Updates the leadlag row on the "Comparison" sheet to match the
#     new totals.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    # Forces a literal text value even when the text looks like a
    # number/date/percentage (e.g. "72.3%", "2026-02-16", "21:45:17"),
    # matching the workbook's inlineStr-style cells instead of letting
    # Excel auto-convert them into numeric/date/percent values. The
    # leading apostrophe is Excel's normal "treat as text" entry marker
    # and is not stored as part of the value; resetting the style back
    # to Normal afterwards drops the quote-prefix formatting flag so the
    # cell keeps the plain default style, matching the source workbook.
    param($range, [string]$text)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("C2").Value = 101
Set-TextValue $wsSummary.Range("D2") "72.3%"
Set-TextValue $wsSummary.Range("E2") "+31.7748%"
Set-TextValue $wsSummary.Range("F2") "+0.3146%"

$wsSummary.Range("C3").Value = 95
Set-TextValue $wsSummary.Range("D3") "52.6%"
Set-TextValue $wsSummary.Range("E3") "+17.9642%"
Set-TextValue $wsSummary.Range("F3") "+0.1891%"

# ---------------------------------------------------------------
# 2. leadlag sheet
# ---------------------------------------------------------------
$wsLead = $wb.Worksheets.Item("leadlag")

# Row 76 - trade #100, DOWN, now CLOSED via time_exit_5min
$wsLead.Range("G76").Value = 68340.813182
Set-TextValue $wsLead.Range("H76") "CLOSED"
$wsLead.Range("I76").Value = -0.0941
$wsLead.Range("J76").Value = -0.9399999999999999
Set-TextValue $wsLead.Range("M76") "time_exit_5min"
$wsLead.Range("N76").Value = 5

# Row 77 - trade #101, UP, now CLOSED via time_exit_5min
$wsLead.Range("G77").Value = 68397.0696
Set-TextValue $wsLead.Range("H77") "CLOSED"
$wsLead.Range("I77").Value = 0.0505
$wsLead.Range("J77").Value = 0.51
Set-TextValue $wsLead.Range("M77") "time_exit_5min"
$wsLead.Range("N77").Value = 5

# New row 97 - trade #122, DOWN, still OPEN
$wsLead.Range("A97").Value = 122
Set-TextValue $wsLead.Range("B97") "2026-02-16"
Set-TextValue $wsLead.Range("C97") "21:45:17"
Set-TextValue $wsLead.Range("D97") "leadlag"
Set-TextValue $wsLead.Range("E97") "DOWN"
$wsLead.Range("F97").Value = 68351.685
Set-TextValue $wsLead.Range("G97") ""
Set-TextValue $wsLead.Range("H97") "OPEN"
$wsLead.Range("I97").Value = 0
$wsLead.Range("J97").Value = 0
$wsLead.Range("K97").Value = 0.75
Set-TextValue $wsLead.Range("L97") "Coinbase leading with -0.107% move"
Set-TextValue $wsLead.Range("M97") ""
$wsLead.Range("N97").Value = 0

# ---------------------------------------------------------------
# 3. All Trades sheet - append the two newly-closed trades
# ---------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# New row 101 - trade #100 (mirrors leadlag row 76 after close)
$wsAll.Range("A101").Value = 100
Set-TextValue $wsAll.Range("B101") "2026-02-16"
Set-TextValue $wsAll.Range("C101") "21:40:03"
Set-TextValue $wsAll.Range("D101") "leadlag"
Set-TextValue $wsAll.Range("E101") "DOWN"
$wsAll.Range("F101").Value = 68276.55
$wsAll.Range("G101").Value = 68340.813182
Set-TextValue $wsAll.Range("H101") "CLOSED"
$wsAll.Range("I101").Value = -0.0941
$wsAll.Range("J101").Value = -0.9399999999999999
$wsAll.Range("K101").Value = 0.75
Set-TextValue $wsAll.Range("L101") "Binance leading with -0.081% move"
Set-TextValue $wsAll.Range("M101") "time_exit_5min"
$wsAll.Range("N101").Value = 5

# New row 102 - trade #101 (mirrors leadlag row 77 after close)
$wsAll.Range("A102").Value = 101
Set-TextValue $wsAll.Range("B102") "2026-02-16"
Set-TextValue $wsAll.Range("C102") "21:40:15"
Set-TextValue $wsAll.Range("D102") "leadlag"
Set-TextValue $wsAll.Range("E102") "UP"
$wsAll.Range("F102").Value = 68362.545
$wsAll.Range("G102").Value = 68397.0696
Set-TextValue $wsAll.Range("H102") "CLOSED"
$wsAll.Range("I102").Value = 0.0505
$wsAll.Range("J102").Value = 0.51
$wsAll.Range("K102").Value = 0.75
Set-TextValue $wsAll.Range("L102") "Binance leading with 0.130% move"
Set-TextValue $wsAll.Range("M102") "time_exit_5min"
$wsAll.Range("N102").Value = 5

# ---------------------------------------------------------------
# 4. Comparison sheet - refreshed leadlag roll-up row
# ---------------------------------------------------------------
$wsComp = $wb.Worksheets.Item("Comparison")

$wsComp.Range("B2").Value = 95
Set-TextValue $wsComp.Range("C2") "52.6%"
Set-TextValue $wsComp.Range("D2") "3.32"
Set-TextValue $wsComp.Range("E2") "+0.5142%"
Set-TextValue $wsComp.Range("F2") "-0.2979%"
Set-TextValue $wsComp.Range("G2") "1.73"
